# Campaign.xlsx update — 6/7/2018 upload
# Update existing campaign rows 2-5 with the latest data, then insert the
# newly-added campaign rows (9 more rows) right after them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the four existing campaign rows with current data
$ws.Range("A2").Value = "apple»SUCCESS»Test»TEST»TEST»TEST»TEST"
$ws.Range("A3").Value = "apple»carrot»robin»TEST»TEST»TEST»TEST"
$ws.Range("A4").Value = "grape»artichoke»sparrow»bear»TEST»TEST"
$ws.Range("A5").Value = "banana»kale»hawk»antelope»thor»TEST»TEST"

# Make room for the 9 new campaign rows right after the existing ones
$ws.Range("A6:A14").Insert()

$ws.Range("A6").Value = "apple»carrot»robin»hamster»spider-man»TEST»TEST"
$ws.Range("A7").Value = "apple»carrot»robin»hamster»spider-man»batman»TEST"
$ws.Range("A8").Value = "banana»kale»hawk»antelope»thor»superman»TEST"
$ws.Range("A9").Value = "grape»artichoke»sparrow»bear»black panther»TEST"
$ws.Range("A10").Value = "TEST Five"
$ws.Range("A11").Value = "TEST Four"
$ws.Range("A12").Value = "TEST TWO"
$ws.Range("A13").Value = "TEST ONE"
$ws.Range("A14").Value = "TEST THREE"
